$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 data")
$ws.Range("B163:B194").Value = 34
